$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry represents one cell whose content changed in the refreshed
# cryptocurrency price/volume snapshot. "ForceText" marks values that look
# like numbers (e.g. "243.09") so we explicitly keep them as text, matching
# how the source data is stored (prices/volumes are text, not numeric cells).
$updates = @(
    @{ Cell = "D2"; Value = "29.087.70"; ForceText = $False },
    @{ Cell = "E2"; Value = "  +0.05%  "; ForceText = $False },
    @{ Cell = "D3"; Value = "1.830.50"; ForceText = $False },
    @{ Cell = "E3"; Value = "  -0.27%  "; ForceText = $False },
    @{ Cell = "D4"; Value = "0.9993"; ForceText = $True },
    @{ Cell = "E4"; Value = "  -0.11%  "; ForceText = $False },
    @{ Cell = "D5"; Value = "243.09"; ForceText = $True },
    @{ Cell = "E5"; Value = "  +0.00%  "; ForceText = $False },
    @{ Cell = "D6"; Value = "0.6280"; ForceText = $True },
    @{ Cell = "E6"; Value = "  +0.30%  "; ForceText = $False },
    @{ Cell = "D7"; Value = "1.001"; ForceText = $True },
    @{ Cell = "E7"; Value = "  -0.06%  "; ForceText = $False },
    @{ Cell = "D8"; Value = "0.07520"; ForceText = $True },
    @{ Cell = "E8"; Value = "  -0.67%  "; ForceText = $False },
    @{ Cell = "D9"; Value = "0.2921"; ForceText = $True },
    @{ Cell = "E9"; Value = "  -0.14%  "; ForceText = $False },
    @{ Cell = "D10"; Value = "23.18"; ForceText = $True },
    @{ Cell = "E10"; Value = "  +2.98%  "; ForceText = $False },
    @{ Cell = "D11"; Value = "0.07672"; ForceText = $True },
    @{ Cell = "E11"; Value = "  -0.79%  "; ForceText = $False },
    @{ Cell = "D12"; Value = "1.830.22"; ForceText = $False },
    @{ Cell = "E12"; Value = "  -0.07%  "; ForceText = $False },
    @{ Cell = "D13"; Value = "5.005"; ForceText = $True },
    @{ Cell = "E13"; Value = "  +0.95%  "; ForceText = $False },
    @{ Cell = "D14"; Value = "0.6677"; ForceText = $True },
    @{ Cell = "E14"; Value = "  +0.59%  "; ForceText = $False },
    @{ Cell = "E15"; Value = "  +0.00%  "; ForceText = $False },
    @{ Cell = "D16"; Value = "0.000009389"; ForceText = $True },
    @{ Cell = "E16"; Value = "  -6.48%  "; ForceText = $False },
    @{ Cell = "E17"; Value = "  -1.14%  "; ForceText = $False },
    @{ Cell = "D18"; Value = "29.069.38"; ForceText = $False },
    @{ Cell = "E18"; Value = "  +0.17%  "; ForceText = $False },
    @{ Cell = "D19"; Value = "2.077.42"; ForceText = $False },
    @{ Cell = "E19"; Value = "  +0.29%  "; ForceText = $False },
    @{ Cell = "E20"; Value = "  +1.70%  "; ForceText = $False },
    @{ Cell = "D21"; Value = "223.02"; ForceText = $True },
    @{ Cell = "E21"; Value = "  -1.60%  "; ForceText = $False },
    @{ Cell = "D22"; Value = "1.002"; ForceText = $True },
    @{ Cell = "E22"; Value = "  -0.04%  "; ForceText = $False },
    @{ Cell = "D23"; Value = "7.109"; ForceText = $True },
    @{ Cell = "E23"; Value = "  -1.06%  "; ForceText = $False },
    @{ Cell = "D24"; Value = "1.001"; ForceText = $True },
    @{ Cell = "E24"; Value = "  -0.16%  "; ForceText = $False },
    @{ Cell = "D25"; Value = "159.77"; ForceText = $True },
    @{ Cell = "E25"; Value = "  +0.65%  "; ForceText = $False },
    @{ Cell = "D26"; Value = "0.1391"; ForceText = $True },
    @{ Cell = "E26"; Value = "  +1.13%  "; ForceText = $False },
    @{ Cell = "D27"; Value = "8.484"; ForceText = $True },
    @{ Cell = "E27"; Value = "  -0.19%  "; ForceText = $False },
    @{ Cell = "D28"; Value = "17.87"; ForceText = $True },
    @{ Cell = "E28"; Value = "  -0.39%  "; ForceText = $False },
    @{ Cell = "D29"; Value = "1.494"; ForceText = $True },
    @{ Cell = "E29"; Value = "  +0.11%  "; ForceText = $False },
    @{ Cell = "D30"; Value = "0.05778"; ForceText = $True },
    @{ Cell = "E30"; Value = "  +10.13%  "; ForceText = $False },
    @{ Cell = "D31"; Value = "4.148"; ForceText = $True },
    @{ Cell = "E31"; Value = "  +1.18%  "; ForceText = $False },
    @{ Cell = "D32"; Value = "4.088"; ForceText = $True },
    @{ Cell = "E32"; Value = "  +1.74%  "; ForceText = $False },
    @{ Cell = "E33"; Value = "  +1.16%  "; ForceText = $False },
    @{ Cell = "B34"; Value = "LidoDAOToken"; ForceText = $False },
    @{ Cell = "C34"; Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"; ForceText = $False },
    @{ Cell = "D34"; Value = "1.830"; ForceText = $True },
    @{ Cell = "E34"; Value = "  -0.57%  "; ForceText = $False },
    @{ Cell = "B35"; Value = "ImmutableX"; ForceText = $False },
    @{ Cell = "C35"; Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; ForceText = $False },
    @{ Cell = "D35"; Value = "0.7383"; ForceText = $True },
    @{ Cell = "E35"; Value = "  +0.67%  "; ForceText = $False },
    @{ Cell = "D36"; Value = "1.136"; ForceText = $True },
    @{ Cell = "E36"; Value = "  -0.16%  "; ForceText = $False },
    @{ Cell = "D37"; Value = "2.669"; ForceText = $True },
    @{ Cell = "E37"; Value = "  -1.11%  "; ForceText = $False },
    @{ Cell = "D38"; Value = "2.766"; ForceText = $True },
    @{ Cell = "E38"; Value = "  +0.10%  "; ForceText = $False },
    @{ Cell = "D39"; Value = "1.216.53"; ForceText = $False },
    @{ Cell = "E39"; Value = "  -1.58%  "; ForceText = $False },
    @{ Cell = "E40"; Value = "  -0.33%  "; ForceText = $False },
    @{ Cell = "D41"; Value = "6.495"; ForceText = $True },
    @{ Cell = "E41"; Value = "  +2.63%  "; ForceText = $False },
    @{ Cell = "D42"; Value = "0.8885"; ForceText = $True },
    @{ Cell = "E42"; Value = "  -1.03%  "; ForceText = $False },
    @{ Cell = "D43"; Value = "1.001"; ForceText = $True },
    @{ Cell = "E43"; Value = "  -0.03%  "; ForceText = $False },
    @{ Cell = "D44"; Value = "102.08"; ForceText = $True },
    @{ Cell = "E44"; Value = "  +0.01%  "; ForceText = $False },
    @{ Cell = "E45"; Value = "  +0.11%  "; ForceText = $False },
    @{ Cell = "D46"; Value = "65.56"; ForceText = $True },
    @{ Cell = "E46"; Value = "  +2.19%  "; ForceText = $False },
    @{ Cell = "D47"; Value = "0.00000000124"; ForceText = $True },
    @{ Cell = "E47"; Value = "  -0.54%  "; ForceText = $False },
    @{ Cell = "D48"; Value = "0.07736"; ForceText = $True },
    @{ Cell = "E48"; Value = "  +15.23%  "; ForceText = $False },
    @{ Cell = "D49"; Value = "0.5088"; ForceText = $True },
    @{ Cell = "E49"; Value = "  -0.52%  "; ForceText = $False },
    @{ Cell = "D50"; Value = "0.4059"; ForceText = $True },
    @{ Cell = "E50"; Value = "  +0.57%  "; ForceText = $False },
    @{ Cell = "D51"; Value = "8.987"; ForceText = $True },
    @{ Cell = "E51"; Value = "  +1.53%  "; ForceText = $False }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $u.Value
}
